$wb = $excel.ActiveWorkbook

# Add new worksheet right after the existing "data" sheet
$dataSheet = $wb.Worksheets.Item("data")
$newSheet = $wb.Worksheets.Add([Type]::Missing, $dataSheet)
$newSheet.Name = "with separators"

$newSheet.Range("A1").Value = "A | B"
$newSheet.Range("B1").Value = "C ! D"
$newSheet.Range("A2").Value = " "
$newSheet.Range("B2").Value = " "

$newSheet.Activate()
$newSheet.Range("B3").Select() | Out-Null
